$wb = $excel.ActiveWorkbook
$win = $excel.ActiveWindow
$win | Get-Member | Out-String | Write-Output
